$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.689.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.752.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '621.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.750.56'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.532'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.490'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000260'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.362.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.751.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.755.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.124'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '506.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.729'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '87.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000138'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +25.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.12'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '45.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '426.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.006.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.42'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.56%  '
